# Apply "average with safety stocks" changes.

$wb = $excel.ActiveWorkbook

# --- Sheet "Productdata": scale InventoryCosts (D), BackorderCosts (F)
#     and LostSale (I) columns for rows 2-11 by 0.0004 ---
$wsProd = $wb.Worksheets.Item("Productdata")

$prodValues = @{
    2  = @{ D = 0.0016; F = 0.0032; I = 0.032  }
    3  = @{ D = 0.0028; F = 0.0056; I = 0.056  }
    4  = @{ D = 0.0024; F = 0.0048; I = 0.048  }
    5  = @{ D = 0.0012; F = 0.0024; I = 0.024  }
    6  = @{ D = 0.0012; F = 0.0024; I = 0.024  }
    7  = @{ D = 0.0012; F = 0.0024; I = 0.024  }
    8  = @{ D = 0.0008; F = 0.0016; I = 0.016  }
    9  = @{ D = 0.0004; F = 0.0008; I = 0.008  }
    10 = @{ D = 0.0004; F = 0.0008; I = 0.008  }
    11 = @{ D = 0.0004; F = 0.0008; I = 0.008  }
}

foreach ($row in $prodValues.Keys) {
    $vals = $prodValues[$row]
    $wsProd.Range("D$row").Value = $vals.D
    $wsProd.Range("F$row").Value = $vals.F
    $wsProd.Range("I$row").Value = $vals.I
}

# --- Sheet "ForcastedStandardDeviation": zero out B:E for rows 9-11 ---
$wsStd = $wb.Worksheets.Item("ForcastedStandardDeviation")

foreach ($row in 9..11) {
    $wsStd.Range("B$row`:E$row").Value = 0
}
